$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (F1:H1), values first
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the formatting from the existing header style (E1) onto the new headers
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

# New boolean columns F:H for rows 2-8 (Outlier flags, MAD-based)
$ws.Range("F2:H8").Value = $false
$ws.Range("G4").Value = $true
